# software-STM: Add function to decode Alarm instruction
#
# The "Alarm" instruction table on Sheet1 documents each byte's meaning for
# the get-alarm-* commands. This change finishes documenting the "get alarm
# extra params" / "set alarm extra params" rows (columns H/I, which record
# whether the day/date-select nibble is "enabled" resp. whether the pattern
# byte has been "fetched") and adds a small legend (M5:N6) explaining the
# "-" / "+" shorthand used throughout the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- fix up existing table cells -------------------------------------------------

# row 9 ("set alarm time params"): enable flag doesn't apply to this row
$ws.Range("K9").Value = "-"

# row 10 ("set alarm extra params"): ringtone/pattern byte decode
$ws.Range("H10").Value = "enable"
$ws.Range("I10").Value = "fetched"

# row 13 ("get alarm extra params"): same decode as row 10
$ws.Range("H13").Value = "enable"
$ws.Range("I13").Value = "fetched"

# --- add the legend at M5:N6 ------------------------------------------------------

$ws.Range("M5").Value = "-"
$ws.Range("M5").HorizontalAlignment = -4108

$ws.Range("N5").Value = "don't care"

$ws.Range("M6").Value = "+"
$ws.Range("M6").HorizontalAlignment = -4108

$ws.Range("N6").Value = "data"
$ws.Range("N6").HorizontalAlignment = -4108

# --- match the author's final selection -------------------------------------------

$ws.Range("A14").Select() | Out-Null
